# Calibrate RCAPEX diet and 18M DSRA funding
#
# 1. CFADS_Calc: fund RCAPEX (col H) for years 6-15 (rows 7-16) and reduce
#    CFADS (col I) accordingly. The Totals row (17) is converted from
#    SUM() formulas into a pasted snapshot of values.
# 2. Ratios: DSCR (col D) is recomputed off the new CFADS figures, the
#    Status column (F) flips for several years, and the post-tenor rows
#    (13-16) lose their (blank) DSCR cell entirely.
# 3. Comparison: the Python/Excel check row is recalibrated and the
#    Excel_Value column (C) is pasted as a static value instead of a
#    live cross-sheet formula.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) CFADS_Calc
# ---------------------------------------------------------------------
$cfads = $wb.Worksheets.Item("CFADS_Calc")

$cfadsRows = @{
    7  = @{ H = 1.2;  I = 16.8  }
    8  = @{ H = 1.5;  I = 20.6  }
    9  = @{ H = 1.65; I = 22.45 }
    10 = @{ H = 1.7;  I = 23.6  }
    11 = @{ H = 1.55; I = 21.25 }
    12 = @{ H = 3.4;  I = 17.6  }
    13 = @{ H = 2.25; I = 17.45 }
    14 = @{ H = 2.0;  I = 15.4  }
    15 = @{ H = 1.5;  I = 13.7  }
    16 = @{ H = 1.25; I = 11.75 }
}

foreach ($r in $cfadsRows.Keys) {
    $vals = $cfadsRows[$r]
    $cfads.Range("H$r").Value = $vals.H
    $cfads.Range("I$r").Value = $vals.I
}

# Totals row: replace the SUM() formulas with a static paste of the
# computed totals.
$cfads.Range("B17").Value = 360
$cfads.Range("C17").Value = 121.7
$cfads.Range("D17").Value = 7.899999999999999
$cfads.Range("E17").Value = 5.400000000000001
$cfads.Range("F17").Value = 10.5
$cfads.Range("G17").Value = 0
$cfads.Range("H17").Value = 36
$cfads.Range("I17").Value = 178.5

# ---------------------------------------------------------------------
# 2) Ratios
# ---------------------------------------------------------------------
$ratios = $wb.Worksheets.Item("Ratios")

# DSCR column (D) precision re-round / recalibration for every year.
$ratiosD = @{
    2  = -0.112994350282
    3  = -0.225988700565
    4  = 0.056497175141
    5  = 0.941619585687
    6  = 1.449999853831
    7  = 1.353333235593
    8  = 1.35158364925
    9  = 1.350726113055
    10 = 1.352569159268
    11 = 1.351425518615
    12 = 1.6440141445
}
foreach ($r in $ratiosD.Keys) {
    $ratios.Range("D$r").Value = $ratiosD[$r]
}

# CFADS column (B) picks up the new RCAPEX-adjusted figures for years 6-15.
$ratiosB = @{
    7  = 16.8
    8  = 20.6
    9  = 22.45
    10 = 23.6
    11 = 21.25
    12 = 17.6
    13 = 17.45
    14 = 15.4
    15 = 13.7
    16 = 11.75
}
foreach ($r in $ratiosB.Keys) {
    $ratios.Range("B$r").Value = $ratiosB[$r]
}

# Post-tenor years (13-16) no longer carry a (blank) DSCR value at all.
$ratios.Range("D13").ClearContents()
$ratios.Range("D14").ClearContents()
$ratios.Range("D15").ClearContents()
$ratios.Range("D16").ClearContents()

# Status column (F) flips for the recalibrated DSCR values.
$ratios.Range("F6").Value  = "PASS"
$ratios.Range("F8").Value  = "PASS"
$ratios.Range("F9").Value  = "PASS"
$ratios.Range("F10").Value = "PASS"
$ratios.Range("F11").Value = "BREACH"
$ratios.Range("F13").Value = "N/A"
$ratios.Range("F14").Value = "N/A"
$ratios.Range("F15").Value = "N/A"
$ratios.Range("F16").Value = "N/A"

# ---------------------------------------------------------------------
# 3) Comparison
# ---------------------------------------------------------------------
$cmp = $wb.Worksheets.Item("Comparison")

$cmp.Range("B2").Value = 196.5
$cmp.Range("C2").Value = 196.5

$cmp.Range("B3").Value = 1.449999853831
$cmp.Range("C3").Value = 1.449999853831

$cmp.Range("B4").Value = 1.6440141445
$cmp.Range("C4").Value = 1.6440141445
